$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) to the table, copying formatting from the
# existing 2019 column (P) so the new cells keep the same styles.

# Header row (row 4): year label 2020, same style as P4.
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null
$ws.Range("Q4").Value = 2020

# Data row (row 5): growth rate value 90.6, same style as P5.
$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5").Value = 90.6

# Restore the active selection to P12, matching the saved view state.
$ws.Range("P12").Select() | Out-Null
